# switched left & right again
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-18 (Plane column C) were "left" -> should become "right"
$ws.Range("C2:C18").Value = "right"

# Rows 36-52 (Plane column C) were "right" -> should become "left"
$ws.Range("C36:C52").Value = "left"

# Update the saved view state (scroll position / selection) to match the edit
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Range("N34").Select()
